# Apply expertise (column G) rating updates on the "10 Rueben Dagenhart" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10 Rueben Dagenhart")
$ws.Activate()

# New ratings for each project row (row number => rating letter, "" clears the cell)
$updates = @{
    2  = "L"
    4  = "M"
    5  = "L"
    6  = "L"
    9  = "L"
    10 = "L"
    12 = "L"
    13 = "M"
    14 = ""
    15 = ""
    16 = ""
    17 = ""
    18 = "L"
    19 = ""
    20 = "L"
    22 = ""
    23 = "M"
    26 = "L"
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $cell = $ws.Cells.Item($row, 7)
    if ($value -eq "") {
        $cell.ClearContents()
    } else {
        $cell.Value = $value
    }
}

# Restore the (non-frozen-pane) selection state recorded in the saved view.
[void]$ws.Range("A10:K10").Select()

[void]$wb.Save()
